# ADD CCR Default BIB
#
# Inserts a new "COUNTERPARTY_BIB_IND_17" / "BE_IND_17" row into the
# "r AnalysisUnit_Variable" sheet (row 16), shifting all subsequent rows
# down by one, and makes that sheet the active tab/selection.

$wb = $excel.ActiveWorkbook

$wsVariable = $wb.Worksheets.Item("r AnalysisUnit_Variable")

# Insert a new blank row at position 16 - this shifts the existing row 16
# (COUNTERPARTY_BIB_IND_16 / BE_IND_16) and everything below it down by one.
$wsVariable.Rows.Item(16).Insert()

# Populate the new row 16 with the CCR default BIB indicator.
$wsVariable.Range("A16").Value = "CREATE/MODIFY"
$wsVariable.Range("B16").Value = "COUNTERPARTY_BIB_IND_17"
$wsVariable.Range("C16").Value = "COUNTERPARTY_BIB_IND_17"
$wsVariable.Range("E16").Value = "COUNTERPARTY_BIB"
$wsVariable.Range("F16").Value = "BE_IND_17"

# Make "r AnalysisUnit_Variable" the active sheet/tab with F16 selected.
$wsVariable.Activate()
$wsVariable.Range("F16").Select()
